# Applies updated JudgeBot opinion text (column C) and, where the
# reclassification changed the recorded outcome, the updated function
# call label (column D) -- "classified every new config".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = @"
MSG: None

MSG: I have recorded the decision to acquire the rights for "Oppenheimer."

"@

$ws.Range("C3").Value = @"
MSG: None

MSG: The decision process did not lead to a consensus, and thus there is no chosen movie for Friday.

"@

$ws.Range("C4").Value = @"
MSG: None

MSG: The decision has been recorded: the movie "Barbie" will be acquired for Friday's showing.

"@

$ws.Range("C5").Value = @"
MSG: None

MSG: The decision has been recorded as no decision being made regarding the movie for Friday.

"@

$ws.Range("C6").Value = @"
MSG: None

MSG: I have recorded the decision to acquire the rights for "Oppenheimer," which will be shown on Friday.

"@

$ws.Range("C7").Value = @"
MSG: None

MSG: I have recorded the decision to acquire the rights for "Oppenheimer" as the movie to be shown on Friday.

"@

$ws.Range("C8").Value = @"
MSG: None

MSG: The decision regarding which movie to show on Friday was not made.

"@

$ws.Range("C9").Value = @"
MSG: None

MSG: The decision has been recorded, and no movie has been selected for Friday.

"@

$ws.Range("C10").Value = @"
MSG: None

MSG: The decision has been made to acquire the rights to "Barbie" for the movie shown on Friday.

"@

$ws.Range("C11").Value = @"
MSG: None

MSG: The decision has been noted as no decision regarding Friday’s movie was reached.

"@

$ws.Range("C12").Value = @"
MSG: None

MSG: The decision to acquire the rights for both movies has been successfully recorded.

"@

$ws.Range("C13").Value = @"
MSG: None

MSG: The decision-making process did not result in a definitive choice regarding the movie to be shown on Friday. Therefore, the outcome is recorded as no decision being made.

"@

$ws.Range("C14").Value = @"
MSG: None

MSG: The decision has been recorded, indicating that no movie was selected for Friday.

"@

$ws.Range("C15").Value = @"
MSG: None

MSG: The committee did not make a decision about which movie to show on Friday.

"@

$ws.Range("C16").Value = @"
MSG: None

MSG: The decision has been recorded as no movie selected for Friday.

"@

$ws.Range("C17").Value = @"
MSG: None

MSG: The decision has been made to acquire the rights for "Barbie" to be shown on Friday.

"@

$ws.Range("C19").Value = @"
MSG: None

MSG: The decision has been recorded: no movie was selected for Friday.

"@

$ws.Range("C20").Value = @"
MSG: None

MSG: The decision has been made to acquire the rights for "Barbie."

"@

$ws.Range("C21").Value = @"
MSG: None

MSG: The decision arrived at by the committee is that no choice of movie for Friday was made.

"@

$ws.Range("C22").Value = @"
MSG: None

MSG: The decision resulted in no movie being selected for Friday.

"@

$ws.Range("C23").Value = @"
MSG: None

MSG: The decision has been recorded to acquire the rights for "Barbie."

"@

$ws.Range("C24").Value = @"
MSG: None

MSG: The rights to both movies have been acquired for the upcoming screening on Friday.

"@

$ws.Range("C25").Value = @"
MSG: None

MSG: The decision has been recorded, and no movie will be acquired for Friday.

"@

$ws.Range("C26").Value = @"
MSG: None

MSG: The decision resulted in no consensus on the movie for Friday, indicating that no choice was made during the discussion.

"@

$ws.Range("C27").Value = @"
MSG: None

MSG: The committee did not reach a decision about which movie to show on Friday.

"@

$ws.Range("C28").Value = @"
MSG: None

MSG: No decision was made regarding the movie.

"@

$ws.Range("C29").Value = @"
MSG: None

MSG: The decision to acquire the rights for both movies, "Barbie" and "Oppenheimer," has been confirmed.

"@

$ws.Range("C30").Value = @"
MSG: None

MSG: The decision has been recorded, and the rights to both movies, "Oppenheimer" and "Barbie," have been acquired for the presentation.

"@

$ws.Range("C31").Value = @"
MSG: None

MSG: The decision about which movie to show on Friday ended without a clear selection. Therefore, I have recorded that there is no decision regarding the movie at this time.

"@

$ws.Range("C32").Value = @"
MSG: None

MSG: The rights to both movies "Oppenheimer" and "Barbie" have been successfully acquired.

"@

$ws.Range("C33").Value = @"
MSG: None

MSG: The rights for both movies have been acquired for the upcoming showing.

"@

$ws.Range("C34").Value = @"
MSG: None

MSG: The rights to both movies will be acquired as the committee expressed a clear intention to show both films.

"@

$ws.Range("C35").Value = @"
MSG: None

MSG: The decision to acquire the rights for "Barbie" has been recorded.

"@

$ws.Range("C36").Value = @"
MSG: None

MSG: The decision about which movie to show on Friday was not reached, so no selection was made.

"@

$ws.Range("C37").Value = @"
MSG: None

MSG: The decision to acquire the rights for both movies has been made successfully.

"@

$ws.Range("C38").Value = @"
MSG: None

MSG: The decision process has ended without a clear choice for Friday's movie, so the conclusion is that no decision can be made at this time.

"@

$ws.Range("C39").Value = @"
MSG: None

MSG: The decision has been made to acquire the rights for "Oppenheimer."

"@

$ws.Range("C40").Value = @"
MSG: None

MSG: The decision has been recorded, and the rights for "Barbie" have been acquired for the assembly on Friday.

"@

$ws.Range("C41").Value = @"
MSG: None

MSG: The decision has been made to acquire the rights for both movies.

"@

$ws.Range("C43").Value = @"
MSG: None

MSG: The rights to both movies have been successfully acquired.

"@

$ws.Range("C44").Value = @"
MSG: None

MSG: The decision has been made to acquire the rights for "Oppenheimer."

"@

$ws.Range("C45").Value = @"
MSG: None

MSG: The rights to both movies have been acquired successfully.

"@

$ws.Range("C46").Value = @"
MSG: None

MSG: The decision has been recorded as no movie selected for Friday.

"@

$ws.Range("C47").Value = @"
MSG: None

MSG: The decision has been recorded, and the rights to "Oppenheimer" will be acquired for Friday's showing.

"@

$ws.Range("C48").Value = @"
MSG: None

MSG: The decision has been made not to select any movie for Friday, as there was no agreement reached between the committee members.

"@

$ws.Range("C50").Value = @"
MSG: None

MSG: The function for no decision has been executed, indicating that no movie was conclusively chosen during the meeting.

"@

$ws.Range("C51").Value = @"
MSG: None

MSG: The decision has been recorded as "no decision," indicating that the committee did not reach a consensus regarding which movie to show on Friday.

"@

$ws.Range("C52").Value = @"
MSG: None

MSG: The decision has been recorded, and no movie will be shown on Friday.

"@

$ws.Range("C53").Value = @"
MSG: None

MSG: The decision to acquire the rights for "Barbie" has been successfully recorded.

"@

$ws.Range("C54").Value = @"
MSG: None

MSG: The decision has been recorded, and the rights for "Barbie" have been acquired.

"@

$ws.Range("C55").Value = @"
MSG: None

MSG: The decision has been recorded as no movie being selected.

"@

$ws.Range("C56").Value = @"
MSG: None

MSG: The decision from the discussion is that no movie was selected for Friday.

"@

$ws.Range("C57").Value = @"
MSG: None

MSG: The decision to show a movie on Friday could not be finalized, leading to no selection being made.

"@

$ws.Range("C58").Value = @"
MSG: None

MSG: The decision has been recorded to acquire the rights for "Oppenheimer."

"@

$ws.Range("C59").Value = @"
MSG: None

MSG: The decision has been logged as "no decision," indicating that no movie was selected for Friday.

"@

$ws.Range("D12").Value = @"
both_movies, 
"@

$ws.Range("D20").Value = @"
Barbie_was_selected, 
"@

$ws.Range("D23").Value = @"
Barbie_was_selected, 
"@

$ws.Range("D30").Value = @"
both_movies, 
"@

$ws.Range("D31").Value = @"
no_decision, 
"@

$ws.Range("D32").Value = @"
both_movies, 
"@

$ws.Range("D43").Value = @"
both_movies, 
"@

$ws.Range("D50").Value = @"
no_decision, 
"@

$ws.Range("D51").Value = @"
no_decision, 
"@

